# Update "想去人数" (want-to-go count) values across all sheets
# This mirrors a scraper re-run that bumped several crowd-interest counters.
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 814
$ws.Range("F3").Value = 14704
$ws.Range("F4").Value = 14943
$ws.Range("F5").Value = 6042
$ws.Range("F8").Value = 1632
$ws.Range("F9").Value = 489
$ws.Range("F11").Value = 1291
$ws.Range("F12").Value = 1956
$ws.Range("F13").Value = 941
$ws.Range("F16").Value = 608
$ws.Range("F18").Value = 3603
$ws.Range("F21").Value = 2647
$ws.Range("F22").Value = 675
$ws.Range("F25").Value = 1907
$ws.Range("F26").Value = 1133
$ws.Range("F27").Value = 1595
$ws.Range("F28").Value = 341
$ws.Range("F30").Value = 7463
$ws.Range("F31").Value = 5175
$ws.Range("F32").Value = 328
$ws.Range("F34").Value = 723
$ws.Range("F35").Value = 722
$ws.Range("F36").Value = 3394
$ws.Range("F39").Value = 358
$ws.Range("F40").Value = 155
$ws.Range("F41").Value = 115
$ws.Range("F42").Value = 4505
$ws.Range("F43").Value = 714
$ws.Range("F44").Value = 27
$ws.Range("F45").Value = 344

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 10
$ws.Range("F20").Value = 112

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7990
$ws.Range("F4").Value = 1102

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7990
$ws.Range("F3").Value = 814
$ws.Range("F5").Value = 1102
$ws.Range("F6").Value = 14704
$ws.Range("F7").Value = 14704
$ws.Range("F8").Value = 6042
$ws.Range("F12").Value = 1632
$ws.Range("F13").Value = 489
$ws.Range("F14").Value = 1291
$ws.Range("F15").Value = 1956
$ws.Range("F19").Value = 3603
$ws.Range("F21").Value = 675
$ws.Range("F24").Value = 1907
$ws.Range("F25").Value = 10
$ws.Range("F30").Value = 1595
$ws.Range("F32").Value = 341
$ws.Range("F34").Value = 7463
$ws.Range("F35").Value = 5175
$ws.Range("F36").Value = 328
$ws.Range("F37").Value = 723
$ws.Range("F38").Value = 3394
$ws.Range("F41").Value = 358
$ws.Range("F43").Value = 115
$ws.Range("F44").Value = 4505
$ws.Range("F45").Value = 714
$ws.Range("F46").Value = 27
$ws.Range("F47").Value = 344

